$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 38.81818
$ws.Range("I6").Value = 38.81818
$ws.Range("K6").Value = 116.45454
$ws.Range("M6").Value = -4.454539999999994
$ws.Range("H38").Value = 1745.4375
$ws.Range("I38").Value = 167.8
$ws.Range("J38").Value = 4374.8335
$ws.Range("K38").Value = 503.4
$ws.Range("L38").Value = 13124.5005
$ws.Range("M38").Value = -131.4
$ws.Range("N38").Value = -13868.5005
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("N43").ClearContents()
$ws.Range("H62").Value = 12000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 12000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H80").Value = 270.66666
$ws.Range("I80").Value = 165.25
$ws.Range("J80").Value = 391.14285
$ws.Range("K80").Value = 495.75
$ws.Range("L80").Value = 1173.42855
$ws.Range("M80").Value = 502.25
$ws.Range("N80").Value = -3169.42855
$ws.Range("H83").Value = 270.66666
$ws.Range("I83").Value = 165.25
$ws.Range("J83").Value = 391.14285
$ws.Range("K83").Value = 1487.25
$ws.Range("L83").Value = 3520.28565
$ws.Range("M83").Value = 3504.75
$ws.Range("N83").Value = -13504.28565
$ws.Range("H105").Value = 17072.428
$ws.Range("J105").Value = 17072.428
$ws.Range("L105").Value = 17072.428
$ws.Range("N105").Value = -24060.428
$ws.Range("H116").Value = 4506.3335
$ws.Range("J116").Value = 6490
$ws.Range("L116").Value = 6490
$ws.Range("N116").Value = -13374
$ws.Range("H137").Value = 2368.2
$ws.Range("I137").Value = 723.75
$ws.Range("J137").Value = 2966.182
$ws.Range("K137").Value = 2171.25
$ws.Range("L137").Value = 8898.545999999998
$ws.Range("M137").Value = 378.75
$ws.Range("N137").Value = -13998.546
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 9600
$ws.Range("I19").Value = 7000
$ws.Range("K19").Value = 7000
$ws.Range("M19").Value = -6771
$ws.Range("H61").Value = 2375
$ws.Range("I61").Value = 2375
$ws.Range("K61").Value = 2375
$ws.Range("M61").Value = -2163
$ws.Range("H88").Value = 1420.0834
$ws.Range("I88").Value = 1630.75
$ws.Range("J88").Value = 998.75
$ws.Range("K88").Value = 1630.75
$ws.Range("L88").Value = 998.75
$ws.Range("M88").Value = -1224.75
$ws.Range("N88").Value = -1810.75
$ws.Range("H91").Value = 1420.0834
$ws.Range("I91").Value = 1630.75
$ws.Range("J91").Value = 998.75
$ws.Range("K91").Value = 1630.75
$ws.Range("L91").Value = 998.75
$ws.Range("M91").Value = -226.75
$ws.Range("N91").Value = -3806.75
$ws.Range("H132").Value = 1654.6364
$ws.Range("I132").Value = 1690.1
$ws.Range("J132").Value = 1300
$ws.Range("K132").Value = 5070.299999999999
$ws.Range("L132").Value = 3900
$ws.Range("M132").Value = -2540.299999999999
$ws.Range("N132").Value = -8960
$ws.Range("H136").Value = 2375
$ws.Range("I136").Value = 2375
$ws.Range("K136").Value = 7125
$ws.Range("M136").Value = -4575
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 986.6
$ws.Range("I80").Value = 596
$ws.Range("K80").Value = 596
$ws.Range("M80").Value = 402
$ws.Range("H83").Value = 986.6
$ws.Range("I83").Value = 596
$ws.Range("K83").Value = 2980
$ws.Range("M83").Value = 2012
$ws.Range("H94").Value = 266.2
$ws.Range("I94").Value = 305.25
$ws.Range("J94").Value = 110
$ws.Range("K94").Value = 305.25
$ws.Range("L94").Value = 110
$ws.Range("M94").Value = 145.75
$ws.Range("N94").Value = -1012
$ws.Range("H99").Value = 1000000000
$ws.Range("I99").Value = 1000000000
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1000000000
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -999998502
$ws.Range("N99").ClearContents()
$ws.Range("H134").Value = 15648.5
$ws.Range("I134").Value = 998
$ws.Range("K134").Value = 2994
$ws.Range("M134").Value = -459
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 1099.8572
$ws.Range("I33").Value = 1099.8572
$ws.Range("K33").Value = 1099.8572
$ws.Range("M33").Value = -720.8571999999999
$ws.Range("H120").Value = 40000
$ws.Range("J120").Value = 40000
$ws.Range("L120").Value = 40000
$ws.Range("N120").Value = -47258
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 1999.5
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H68").Value = 500
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 500
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 16417019
$ws.Range("J11").Value = 26750012
$ws.Range("L11").Value = 26750012
$ws.Range("N11").Value = -26750290
$ws.Range("H80").Value = 2150.6667
$ws.Range("I80").Value = 1951.125
$ws.Range("J80").Value = 2549.75
$ws.Range("K80").Value = 1951.125
$ws.Range("L80").Value = 2549.75
$ws.Range("M80").Value = -953.125
$ws.Range("N80").Value = -4545.75
$ws.Range("H83").Value = 2150.6667
$ws.Range("I83").Value = 1951.125
$ws.Range("J83").Value = 2549.75
$ws.Range("K83").Value = 9755.625
$ws.Range("L83").Value = 12748.75
$ws.Range("M83").Value = -4763.625
$ws.Range("N83").Value = -22732.75
$ws.Range("H97").Value = 843.3
$ws.Range("I97").Value = 224.8
$ws.Range("J97").Value = 1461.8
$ws.Range("K97").Value = 224.8
$ws.Range("L97").Value = 1461.8
$ws.Range("M97").Value = 271.2
$ws.Range("N97").Value = -2453.8
$ws.Range("H101").Value = 23664.666
$ws.Range("J101").Value = 23664.666
$ws.Range("L101").Value = 23664.666
$ws.Range("N101").Value = -30154.666
$ws.Range("H122").Value = 2934.8572
$ws.Range("I122").Value = 2394
$ws.Range("J122").Value = 3151.2
$ws.Range("K122").Value = 7182
$ws.Range("L122").Value = 9453.599999999999
$ws.Range("M122").Value = -4732
$ws.Range("N122").Value = -14353.6
$ws.Range("H132").Value = 2598.5
$ws.Range("I132").Value = 2598.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7795.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5265.5
$ws.Range("N132").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 17500
$ws.Range("J64").Value = 17500
$ws.Range("L64").Value = 17500
$ws.Range("N64").Value = -17950
$ws.Range("H67").Value = 17500
$ws.Range("J67").Value = 17500
$ws.Range("L67").Value = 17500
$ws.Range("N67").Value = -19060
$ws.Range("H82").Value = 3036.2222
$ws.Range("I82").Value = 469.25
$ws.Range("J82").Value = 5089.8
$ws.Range("K82").Value = 469.25
$ws.Range("L82").Value = 5089.8
$ws.Range("M82").Value = -108.25
$ws.Range("N82").Value = -5811.8
$ws.Range("H85").Value = 3036.2222
$ws.Range("I85").Value = 469.25
$ws.Range("J85").Value = 5089.8
$ws.Range("K85").Value = 469.25
$ws.Range("L85").Value = 5089.8
$ws.Range("M85").Value = 778.75
$ws.Range("N85").Value = -7585.8
$ws.Range("H100").Value = 10000
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 10000
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 10000
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -11082
$ws.Range("H132").Value = 3399.5715
$ws.Range("I132").Value = 3674.25
$ws.Range("J132").Value = 3033.3333
$ws.Range("K132").Value = 11022.75
$ws.Range("L132").Value = 9099.999899999999
$ws.Range("M132").Value = -8492.75
$ws.Range("N132").Value = -14159.9999
$ws.Range("H136").Value = 6002.5
$ws.Range("I136").Value = 6002.5
$ws.Range("K136").Value = 18007.5
$ws.Range("M136").Value = -15457.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 13875
$ws.Range("J31").Value = 13875
$ws.Range("L31").Value = 13875
$ws.Range("N31").Value = -14571
$ws.Range("H96").Value = 868.1667
$ws.Range("I96").Value = 1125
$ws.Range("J96").Value = 739.75
$ws.Range("K96").Value = 1125
$ws.Range("L96").Value = 739.75
$ws.Range("M96").Value = 248
$ws.Range("N96").Value = -3485.75
$ws.Range("H122").Value = 2924.5625
$ws.Range("I122").Value = 2663.4285
$ws.Range("K122").Value = 7990.2855
$ws.Range("M122").Value = -5540.2855
$ws.Range("H132").Value = 1583.75
$ws.Range("I132").Value = 1111.6666
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 3334.9998
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -804.9998000000001
$ws.Range("N132").Value = -14060
